$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new "Wins" / "Losses" / "Ties" columns (AD1:AF1), styled like
# the rest of the header row (bold, bordered, centered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header cells (bold font, thin
# border, centered horizontal/top alignment) by copying the format from
# the last existing header cell (AC1) onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-38: season record (Wins=91, Losses=71, Ties=0) repeated for
# every player row.
$lastRow = 38
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}
